$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 1.360200666666667
$ws.Range("N2").Value = 4.080602
$ws.Range("O2").Value = 0.014316771935026
$ws.Range("P2").Value = 0.014316771935026
$ws.Range("Q2").Value = 0.7441925353464445
$ws.Range("R2").Value = 6.697732818117999
$ws.Range("S2").Value = 0.014316771935026
$ws.Range("T2").Value = 0.014316771935026

# Row 3 updates
$ws.Range("O3").Value = 0.7566069373802504
$ws.Range("P3").Value = 0.7566069373802504
$ws.Range("S3").Value = 0.7566069373802504
$ws.Range("T3").Value = 0.7566069373802504

# Row 4 updates
$ws.Range("M4").Value = 21.763965
$ws.Range("N4").Value = 65.291895
$ws.Range("O4").Value = 0.2290762906847235
$ws.Range("P4").Value = 0.2290762906847235
$ws.Range("Q4").Value = 11.907493276145
$ws.Range("R4").Value = 107.167439485305
$ws.Range("S4").Value = 0.2290762906847235
$ws.Range("T4").Value = 0.2290762906847235
